$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I11").Value = 'sd'
$ws.Range("J11").Value = 'Statement-non-opinion'
$ws.Range("I15").Value = 'sd'
$ws.Range("J15").Value = 'Statement-non-opinion'
$ws.Range("I20").Value = 'aa'
$ws.Range("J20").Value = 'Agree/Accept'
$ws.Range("I21").Value = 'b'
$ws.Range("J21").Value = 'Acknowledge (Backchannel)'
$ws.Range("I28").Value = 'b'
$ws.Range("J28").Value = 'Acknowledge (Backchannel)'
$ws.Range("I29").Value = 'sd'
$ws.Range("J29").Value = 'Statement-non-opinion'
$ws.Range("I31").Value = 'sv'
$ws.Range("J31").Value = 'Statement-opinion'
$ws.Range("I36").Value = 'b'
$ws.Range("J36").Value = 'Acknowledge (Backchannel)'
$ws.Range("I52").Value = 'b'
$ws.Range("J52").Value = 'Acknowledge (Backchannel)'
$ws.Range("I57").Value = '%'
$ws.Range("J57").Value = 'Uninterpretable'
$ws.Range("I71").Value = 'aa'
$ws.Range("J71").Value = 'Agree/Accept'
$ws.Range("I78").Value = 'sd'
$ws.Range("J78").Value = 'Statement-non-opinion'
$ws.Range("I79").Value = 'b'
$ws.Range("J79").Value = 'Acknowledge (Backchannel)'
$ws.Range("I82").Value = 'b'
$ws.Range("J82").Value = 'Acknowledge (Backchannel)'
$ws.Range("I96").Value = 'sd'
$ws.Range("J96").Value = 'Statement-non-opinion'
$ws.Range("I109").Value = 'ba'
$ws.Range("J109").Value = 'Appreciation'
$ws.Range("I118").Value = 'b'
$ws.Range("J118").Value = 'Acknowledge (Backchannel)'
$ws.Range("I130").Value = 'ba'
$ws.Range("J130").Value = 'Appreciation'
$ws.Range("I134").Value = 'sd'
$ws.Range("J134").Value = 'Statement-non-opinion'
$ws.Range("I135").Value = 'sd'
$ws.Range("J135").Value = 'Statement-non-opinion'
$ws.Range("I137").Value = 'sd'
$ws.Range("J137").Value = 'Statement-non-opinion'
$ws.Range("I138").Value = 'sv'
$ws.Range("J138").Value = 'Statement-opinion'
$ws.Range("I144").Value = 'ba'
$ws.Range("J144").Value = 'Appreciation'
$ws.Range("I158").Value = 'ba'
$ws.Range("J158").Value = 'Appreciation'
$ws.Range("I159").Value = 'b'
$ws.Range("J159").Value = 'Acknowledge (Backchannel)'
$ws.Range("I161").Value = 'b'
$ws.Range("J161").Value = 'Acknowledge (Backchannel)'
$ws.Range("I170").Value = '%'
$ws.Range("J170").Value = 'Uninterpretable'
$ws.Range("I180").Value = '%'
$ws.Range("J180").Value = 'Uninterpretable'
$ws.Range("I207").Value = 'ba'
$ws.Range("J207").Value = 'Appreciation'
$ws.Range("I218").Value = 'sv'
$ws.Range("J218").Value = 'Statement-opinion'
$ws.Range("I228").Value = 'sd'
$ws.Range("J228").Value = 'Statement-non-opinion'
$ws.Range("I235").Value = 'sd'
$ws.Range("J235").Value = 'Statement-non-opinion'
$ws.Range("I239").Value = 'sd'
$ws.Range("J239").Value = 'Statement-non-opinion'
$ws.Range("I245").Value = 'sv'
$ws.Range("J245").Value = 'Statement-opinion'
$ws.Range("I258").Value = 'ba'
$ws.Range("J258").Value = 'Appreciation'
$ws.Range("I259").Value = 'sd'
$ws.Range("J259").Value = 'Statement-non-opinion'
$ws.Range("I263").Value = 'aa'
$ws.Range("J263").Value = 'Agree/Accept'
$ws.Range("I296").Value = 'ba'
$ws.Range("J296").Value = 'Appreciation'
